# Generate Report for Handoff
# - Flip status from "Handed back: in sync with en-US" to "Ready for handoff"
# - Bump the "Latest HO Xliff Generate Date" / handoff datetimes forward
# - Narrow the (now shorter) status columns

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Status column(s): "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Datetime bumps
$overview.Range("G2").Value = "2016-08-19 23:04:30"
$dede.Range("H2").Value = "2016-08-19 23:04:30"
$zhcn.Range("H2").Value = "2016-08-19 23:04:25"

# Column widths shrink now that "Ready for handoff" is shorter than
# "Handed back: in sync with en-US" (~17.2 chars vs ~30 chars wide)
$overview.Range("E1").ColumnWidth = 16.3
$overview.Range("F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
